# Fill the BIK check-report table with the new bank list (rows 2-11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "КУ ОАО АКБ `"УНИВЕРСАЛЬНЫЙ КРЕДИТ`" - ГК `"АСВ`" "
$ws.Range("B2").Value = "г Москва "
$ws.Range("C2").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30101810345250000651 "
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "044525651"

$ws.Range("A3").Value = "КУ `"НАЦКОРПБАНК`" (АО) - ГК `"АСВ`" "
$ws.Range("B3").Value = "г Москва "
$ws.Range("C3").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30101810945250000653 "
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "044525653"

$ws.Range("A4").Value = "КУ ОАО АКБ `"РУССКИЙ ЗЕМЕЛЬНЫЙ БАНК`" -ГК `"АСВ`" "
$ws.Range("B4").Value = "г Москва "
$ws.Range("C4").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "30101810545250000684 "
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "044525684"

$ws.Range("A5").Value = "КУ ООО КБ `"ИНСТРОЙБАНК`" ГК `"АСВ`" "
$ws.Range("B5").Value = "г Москва "
$ws.Range("C5").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "30101810445250000690 "
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "044525690"

$ws.Range("A6").Value = "КУ ОАО АКБ `"ЛЕСБАНК`"-ГК `"АСВ`" "
$ws.Range("B6").Value = "г Москва "
$ws.Range("C6").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "30101810645250000694 "
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "044525694"

$ws.Range("A7").Value = "КУ ПАО `"ТАЙМ БАНК`" - ГК `"АСВ`" "
$ws.Range("B7").Value = "г Москва "
$ws.Range("C7").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "30101810445250000713 "
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "044525713"

$ws.Range("A8").Value = "КУ ООО КБ `"ПЕРВЫЙ ДЕПОЗИТНЫЙ`" - ГК `"АСВ`" "
$ws.Range("B8").Value = "г Москва "
$ws.Range("C8").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30101810845250000737 "
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "044525737"

$ws.Range("A9").Value = "КУ АКБ `"МФТ-БАНК`" (ОАО) - ГК `"АСВ`" "
$ws.Range("B9").Value = "г Москва "
$ws.Range("C9").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "30101810745250000743 "
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "044525743"

$ws.Range("A10").Value = "КУ АО `"СМАРТБАНК`" - ГК `"АСВ`" "
$ws.Range("B10").Value = "г Москва "
$ws.Range("C10").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30101810645250000746 "
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "044525746"

$ws.Range("A11").Value = "КУ ЗАО МЕЖДУНАРОДНЫЙ ПРОМЫШЛЕННЫЙ БАНК-ГК АСВ "
$ws.Range("B11").Value = "г Москва "
$ws.Range("C11").Value = "ул Высоцкого, 4 "
# KorrSchet / BIK look numeric but must stay text (leading zeros, trailing space, > 15 significant digits)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30101810000000000748 "
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "044525748"

